$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the data area first so cells that should become empty are actually empty
$ws.Range("A1:C9").ClearContents()

# Row 1 - headers (unchanged)
$ws.Range("A1").Value = "level_i"
$ws.Range("B1").Value = "level_ii"
$ws.Range("C1").Value = "level_iii"

# Row 2 - student_id -> sid
$ws.Range("A2").Value = "sid"
$ws.Range("B2").Value = "sid"
$ws.Range("C2").Value = "sid"

# Row 3 - name -> sname
$ws.Range("A3").Value = "sname"
$ws.Range("B3").Value = "sname"
$ws.Range("C3").Value = "sname"

# Row 4 - dob (A,B only, C left blank)
$ws.Range("A4").Value = "dob"
$ws.Range("B4").Value = "dob"

# Row 5 - gender (A only)
$ws.Range("A5").Value = "gender"

# Row 6 - raceeth (A only)
$ws.Range("A6").Value = "raceeth"

# Row 7 - teacher_id -> tid (A only)
$ws.Range("A7").Value = "tid"

# Row 8 - teacher_name -> tname (A,B,C)
$ws.Range("A8").Value = "tname"
$ws.Range("B8").Value = "tname"
$ws.Range("C8").Value = "tname"

# Row 9 - zip_code -> zip (A,B only)
$ws.Range("A9").Value = "zip"
$ws.Range("B9").Value = "zip"
